$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph 3 ("That is my task to push by my own") gets a grammar-check
#    style split: the trailing word "own" is pulled into its own run and
#    wrapped with proofErr gramStart/gramEnd markers, while the rest of the
#    sentence becomes a separate run ending in a trailing space.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3).Range

# Preserve the paragraph's own identity attributes (paraId/textId/rsid*)
# instead of hard-coding them, by pulling them off its current OOXML.
$p3Attrs = ""
if ($p3.WordOpenXML -match '<w:p\s+([^>]*)>') {
    $p3Attrs = $Matches[1]
}

$p3Xml = '<?xml version="1.0"?><?mso-application progid="Word.Document"?>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' + `
    'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
    '<w:body><w:p ' + $p3Attrs + '>' + `
    '<w:r><w:t xml:space="preserve">That is my task to push by my </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>own</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '</w:p></w:body></w:document>'
$p3.InsertXML($p3Xml)

# ---------------------------------------------------------------------------
# 2) A brand-new paragraph "My 4th task" is added right after it, with the
#    ordinal suffix "th" raised to superscript.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3).Range
$p3.InsertParagraphAfter()

$p4 = $d.Paragraphs.Item(4).Range
$p4.Text = "My 4th task"

$supRange = $d.Range($p4.Start, $p4.End)
$supRange.Find.Execute("th", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$supRange.Font.Superscript = $true
